$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 28, pushing the existing rows 28-31 down to 29-32.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new record.
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44736
$ws.Range("D28").NumberFormat = $ws.Range("D27").NumberFormat
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112045
$ws.Range("G28").Value = "Zapallo"
$ws.Range("H28").Value = "Camote"
$ws.Range("I28").Value = "1a (guarda)"
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 780
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = 790
$ws.Range("N28").Value = "$/kilo (volumen en unidades)"
$ws.Range("O28").Value = "Región de O'Higgins"
$ws.Range("P28").Value = 790
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"
